$d = $word.ActiveDocument

# Locate the paragraph that ends with the "...land on 10, 100, and 1,000." text
# (this is the paragraph that currently carries the _GoBack bookmark at its end),
# and the following (currently empty) list paragraph.
$paraGoal = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*land on 10, 100, and 1,000.*") {
        $paraGoal = $p
    }
}
$paraEmpty = $paraGoal.Next()

# 1) Remove the _GoBack bookmark from the end of $paraGoal and replace it with a
#    plain run containing a single space.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goalEnd = $paraGoal.Range.End - 1
$goalRange = $d.Range($goalEnd, $goalEnd)
$goalRange.InsertAfter(" ")

# 2) Fill in the previously empty paragraph with the "sub-goals" text.
$emptyRange = $paraEmpty.Range
$emptyRange.End = $emptyRange.End - 1
$emptyRange.Text = "So the sub-goals are finding the fingers on three numbers. The numbers being 10, 100, and 1,000. The only constraint is how the girl counts out the numbers on her fingers. She starts with 1 on her thumb, going 2 on her first finger, 3 on her middle finger, 4 on her ring finger, and five on her pinky. She then reverses it going 6 on her ring finger, 7 on her middle, 8 on her first finger, and 9 on her thumb. This continues until she reaches each of the goals."

# 3) Insert a new list paragraph describing the solutions reached, made up of two
#    separate runs. Build it as two temporary paragraphs (so each gets its own
#    run), then join them back into a single paragraph by deleting the paragraph
#    mark between them - this keeps the two runs distinct instead of merging them.
$p1Range = $paraEmpty.Range
$p1Range.End = $p1Range.End - 1
$p1Range.Collapse(0)
$p1Range.InsertParagraphAfter()
$paraSolA = $paraEmpty.Next()

$solARange = $paraSolA.Range
$solARange.End = $solARange.End - 1
$solARange.Text = "So the solution to the first goal is given in the question. She would reach 10 on her first finger. The ways to figure the other do can either be done by actually counting out the numbers on my own hand or create a chart to find a pattern. With that done the solutions are that the girl would reach 100 on"

$p2Range = $paraSolA.Range
$p2Range.End = $p2Range.End - 1
$p2Range.Collapse(0)
$p2Range.InsertParagraphAfter()
$paraSolB = $paraSolA.Next()

$solBRange = $paraSolB.Range
$solBRange.End = $solBRange.End - 1
$solBRange.Text = " her ring finger and 1,000 on her first finger."

$joinPos = $paraSolA.Range.End - 1
$joinRange = $d.Range($joinPos, $joinPos + 1)
$joinRange.Delete()

# 4) Insert one more new (empty) list paragraph after that one, and give it the
#    _GoBack bookmark that was removed in step 1. A placeholder character is
#    inserted first so the (non-collapsed) bookmark can be created reliably,
#    then the placeholder is deleted, leaving a collapsed bookmark behind.
$p3Range = $paraSolA.Range
$p3Range.End = $p3Range.End - 1
$p3Range.Collapse(0)
$p3Range.InsertParagraphAfter()
$paraFinal = $paraSolA.Next()

$finalStart = $paraFinal.Range.Start
$placeholderRange = $d.Range($finalStart, $finalStart)
$placeholderRange.InsertBefore("X")
$bmRange = $d.Range($finalStart, $finalStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$xRange = $d.Range($finalStart, $finalStart + 1)
$xRange.Delete()
